# Insert a new weekly price record at row 32 ("Florida King" durazno,
# Región de Coquimbo), pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(32).Insert()

$ws.Range("A32").Value = 1
$ws.Range("B32").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C32").Value = "Arica y Parinacota"
$ws.Range("D32").Value = 45210
$ws.Range("E32").Value = 15
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100103
$ws.Range("H32").Value = "Frutos de hueso (carozo)"
$ws.Range("I32").Value = 100103004
$ws.Range("J32").Value = "Durazno"
$ws.Range("K32").Value = "Florida King"
$ws.Range("L32").Value = "Segunda"
$ws.Range("M32").Value = 300
$ws.Range("N32").Value = 29000
$ws.Range("O32").Value = 30000
$ws.Range("P32").Value = 29500
$ws.Range("Q32").Value = "`$/bandeja 10 kilos granel"
$ws.Range("R32").Value = "Región de Coquimbo"
$ws.Range("S32").Value = 2950
$ws.Range("T32").Value = 10
